$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (column E) values for rows 19-30 so that the
# period sequence runs consecutively starting at 2006 (row 19) through
# 2105 (row 30). This reshuffles the previously-reversed 2007-2012 /
# 2101-2105 block and adds the new 2006 entry at row 19, pushing the
# remaining periods down by one and wrapping the old 2006 value (with
# its distinct "Valor Mora" of 28090) onto the final row.
$ws.Range("E19").Value = "2006"
$ws.Range("E20").Value = "2007"
$ws.Range("E21").Value = "2008"
$ws.Range("E22").Value = "2009"
$ws.Range("E23").Value = "2010"
$ws.Range("E24").Value = "2011"
$ws.Range("E25").Value = "2012"
$ws.Range("E26").Value = "2101"
$ws.Range("E27").Value = "2102"
$ws.Range("E28").Value = "2103"
$ws.Range("E29").Value = "2104"
$ws.Range("E30").Value = "2105"

# "Valor Mora" (column F): the distinctive 28090 value follows the 2006
# period, so it moves from row 19 to row 30; row 19 now takes the
# standard 35112 value.
$ws.Range("F19").Value = 35112
$ws.Range("F30").Value = 28090
